$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records (2020-12-30 .. 2021-01-09) appended as rows 309-319
$data = @(
    @("2020-12-30", 2142, 633, 1073, 436, 2803, 353, 504, 1946, 29.55, 50.09, 20.35, 12.59, 17.98, 69.43000000000001),
    @("2020-12-31", 2142, 676, 1043, 423, 2804, 348, 476, 1980, 31.56, 48.69, 19.75, 12.41, 16.98, 70.61),
    @("2021-01-01", 2141, 663, 969, 509, 2803, 342, 497, 1964, 30.97, 45.26, 23.77, 12.2, 17.73, 70.06999999999999),
    @("2021-01-02", 2141, 648, 918, 575, 2801, 326, 505, 1970, 30.27, 42.88, 26.86, 11.64, 18.03, 70.33),
    @("2021-01-03", 2137, 649, 933, 555, 2802, 321, 505, 1976, 30.37, 43.66, 25.97, 11.46, 18.02, 70.52),
    @("2021-01-04", 2140, 591, 1027, 522, 2799, 326, 503, 1970, 27.62, 47.99, 24.39, 11.65, 17.97, 70.38),
    @("2021-01-05", 2141, 586, 1071, 484, 2802, 303, 520, 1979, 27.37, 50.02, 22.61, 10.81, 18.56, 70.63),
    @("2021-01-06", 2156, 593, 1108, 455, 2805, 306, 503, 1996, 27.5, 51.39, 21.1, 10.91, 17.93, 71.16),
    @("2021-01-07", 2146, 584, 1068, 494, 2802, 295, 494, 2013, 27.21, 49.77, 23.02, 10.53, 17.63, 71.84),
    @("2021-01-08", 2146, 579, 1114, 453, 2803, 299, 500, 2004, 26.98, 51.91, 21.11, 10.67, 17.84, 71.48999999999999),
    @("2021-01-09", 2150, 556, 1082, 512, 2800, 291, 469, 2040, 25.86, 50.33, 23.81, 10.39, 16.75, 72.86)
)

$startRow = 309
$endRow = $startRow + $data.Length - 1

# Format column A as text first so the DATE strings (e.g. "2020-12-30")
# are stored as literal text, matching the rest of the DATE column,
# instead of being auto-coerced into Excel date serials.
$dateRangeAddr = "A" + $startRow + ":A" + $endRow
$dateRange = $ws.Range($dateRangeAddr)
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Restore the default (unformatted) style on column A now that the
# values are committed as text, so the new cells match the plain
# (style-less) DATE cells above them.
$dateRange.Style = "Normal"
